$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting for new rows 15-35 (column A date style) by copying row 14 down
$ws.Range("A14:D14").Copy()
$ws.Range("A15:D35").PasteSpecial(-4122)

$ws.Range("A2").Value = "2026-06-02"
$ws.Range("B2").Value = "TV"
$ws.Range("C2").Value = "Spend"
$ws.Range("D2").Value = 119
$ws.Range("A3").Value = "2026-06-09"
$ws.Range("B3").Value = "Radio"
$ws.Range("C3").Value = "Spend"
$ws.Range("D3").Value = 59
$ws.Range("A4").Value = "2026-04-28"
$ws.Range("B4").Value = "TV"
$ws.Range("C4").Value = "Spend"
$ws.Range("D4").Value = 181
$ws.Range("A5").Value = "2026-05-19"
$ws.Range("B5").Value = "TV"
$ws.Range("C5").Value = "Spend"
$ws.Range("D5").Value = 56
$ws.Range("A6").Value = "2026-05-05"
$ws.Range("B6").Value = "TV"
$ws.Range("C6").Value = "GRPs"
$ws.Range("D6").Value = 8
$ws.Range("A7").Value = "2026-05-05"
$ws.Range("B7").Value = "TV"
$ws.Range("C7").Value = "Spend"
$ws.Range("D7").Value = 126
$ws.Range("A8").Value = "2026-05-26"
$ws.Range("B8").Value = "TV"
$ws.Range("C8").Value = "Spend"
$ws.Range("D8").Value = 90
$ws.Range("A9").Value = "2026-05-19"
$ws.Range("B9").Value = "TV"
$ws.Range("C9").Value = "Spend"
$ws.Range("D9").Value = 56
$ws.Range("A10").Value = "2026-06-09"
$ws.Range("B10").Value = "TV"
$ws.Range("C10").Value = "GRPs"
$ws.Range("D10").Value = 8
$ws.Range("A11").Value = "2026-06-16"
$ws.Range("B11").Value = "TV"
$ws.Range("C11").Value = "GRPs"
$ws.Range("D11").Value = 3
$ws.Range("A12").Value = "2026-06-09"
$ws.Range("B12").Value = "TV"
$ws.Range("C12").Value = "Spend"
$ws.Range("D12").Value = 98
$ws.Range("A13").Value = "2026-05-12"
$ws.Range("B13").Value = "TV"
$ws.Range("C13").Value = "GRPs"
$ws.Range("D13").Value = 9
$ws.Range("A14").Value = "2026-05-12"
$ws.Range("B14").Value = "TV"
$ws.Range("C14").Value = "Spend"
$ws.Range("D14").Value = 51
$ws.Range("A15").Value = "2026-05-19"
$ws.Range("B15").Value = "TV"
$ws.Range("C15").Value = "GRPs"
$ws.Range("D15").Value = 10
$ws.Range("A16").Value = "2026-04-28"
$ws.Range("B16").Value = "Radio"
$ws.Range("C16").Value = "Spend"
$ws.Range("D16").Value = 98
$ws.Range("A17").Value = "2026-05-26"
$ws.Range("B17").Value = "Radio"
$ws.Range("C17").Value = "Spend"
$ws.Range("D17").Value = 200
$ws.Range("A18").Value = "2026-04-28"
$ws.Range("B18").Value = "Radio"
$ws.Range("C18").Value = "GRPs"
$ws.Range("D18").Value = 5
$ws.Range("A19").Value = "2026-05-26"
$ws.Range("B19").Value = "TV"
$ws.Range("C19").Value = "GRPs"
$ws.Range("D19").Value = 3
$ws.Range("A20").Value = "2026-05-26"
$ws.Range("B20").Value = "TV"
$ws.Range("C20").Value = "GRPs"
$ws.Range("D20").Value = 3
$ws.Range("A21").Value = "2026-05-12"
$ws.Range("B21").Value = "TV"
$ws.Range("C21").Value = "Spend"
$ws.Range("D21").Value = 51
$ws.Range("A22").Value = "2026-05-05"
$ws.Range("B22").Value = "Radio"
$ws.Range("C22").Value = "Spend"
$ws.Range("D22").Value = 77
$ws.Range("A23").Value = "2026-04-28"
$ws.Range("B23").Value = "TV"
$ws.Range("C23").Value = "GRPs"
$ws.Range("D23").Value = 1
$ws.Range("A24").Value = "2026-05-05"
$ws.Range("B24").Value = "TV"
$ws.Range("C24").Value = "Spend"
$ws.Range("D24").Value = 126
$ws.Range("A25").Value = "2026-05-19"
$ws.Range("B25").Value = "Radio"
$ws.Range("C25").Value = "GRPs"
$ws.Range("D25").Value = 10
$ws.Range("A26").Value = "2026-05-26"
$ws.Range("B26").Value = "Radio"
$ws.Range("C26").Value = "GRPs"
$ws.Range("D26").Value = 1
$ws.Range("A27").Value = "2026-06-02"
$ws.Range("B27").Value = "Radio"
$ws.Range("C27").Value = "Spend"
$ws.Range("D27").Value = 147
$ws.Range("A28").Value = "2026-05-12"
$ws.Range("B28").Value = "TV"
$ws.Range("C28").Value = "GRPs"
$ws.Range("D28").Value = 9
$ws.Range("A29").Value = "2026-06-02"
$ws.Range("B29").Value = "Radio"
$ws.Range("C29").Value = "GRPs"
$ws.Range("D29").Value = 7
$ws.Range("A30").Value = "2026-06-16"
$ws.Range("B30").Value = "Radio"
$ws.Range("C30").Value = "GRPs"
$ws.Range("D30").Value = 1
$ws.Range("A31").Value = "2026-06-16"
$ws.Range("B31").Value = "Radio"
$ws.Range("C31").Value = "Spend"
$ws.Range("D31").Value = 76
$ws.Range("A32").Value = "2026-06-09"
$ws.Range("B32").Value = "Radio"
$ws.Range("C32").Value = "GRPs"
$ws.Range("D32").Value = 5
$ws.Range("A33").Value = "2026-06-16"
$ws.Range("B33").Value = "TV"
$ws.Range("C33").Value = "Spend"
$ws.Range("D33").Value = 144
$ws.Range("A34").Value = "2026-05-19"
$ws.Range("B34").Value = "Radio"
$ws.Range("C34").Value = "GRPs"
$ws.Range("D34").Value = 10
$ws.Range("A35").Value = "2026-06-09"
$ws.Range("B35").Value = "Radio"
$ws.Range("C35").Value = "GRPs"
$ws.Range("D35").Value = 5
